# tdf125573_FontWorkScaleX.pptx: change FontWork run font from Cambria to
# Calibri (both the Latin and East-Asian typeface) on slide 1's WordArt
# shape, which holds the two paragraphs
#   "This is a long, long first line."
#   "Second line short"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$line1 = "This is a long, long first line."
$line2 = "Second line short"

# The two paragraphs currently share identical run formatting (same <a:rPr>
# contents, just different <a:t>). Setting Font.NameFarEast on a sub-range
# (a single paragraph/run) only ever touches the *first* run of the shape in
# this host, no matter which characters are targeted. Work around that by
# temporarily merging the two paragraphs back into a single run (so there is
# only one, unambiguous "first run"), changing the font there, and then
# re-splitting the text back into its original two paragraphs - the
# paragraph mark picks up a duplicate of that (now-updated) run formatting
# for both resulting paragraphs.
$tr.Text = "$line1$line2"

$tr.Font.Name = "Calibri"
$tr.Font.NameFarEast = "Calibri"

$tr.Text = "$line1`r$line2"
